# Applies the Case_0_197 vm_pu.xlsx re-run results ("case with 380 kV done").
# For each data row (2-25), columns B:F and I:N hold per-bus voltage magnitudes
# (column A is the index, G/H are untouched in this run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.047222114997
$rowVals[0,2] = 1.054555962832996
$rowVals[0,3] = 1.054557186606952
$rowVals[0,4] = 1.064901180420221
$ws.Range("B2:F2").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041657943018336
$rowVals[0,1] = 1.052272332560762
$rowVals[0,2] = 1.057298848598588
$rowVals[0,3] = 1.057300069001437
$rowVals[0,4] = 1.06761586801368
$rowVals[0,5] = 1.053766679970322
$ws.Range("I2:N2").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.048186582068708
$rowVals[0,2] = 1.055410390718746
$rowVals[0,3] = 1.055406761641134
$rowVals[0,4] = 1.065816379215227
$ws.Range("B3:F3").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041826882224083
$rowVals[0,1] = 1.052885111579656
$rowVals[0,2] = 1.057966437851084
$rowVals[0,3] = 1.057962818048283
$rowVals[0,4] = 1.068346111644549
$rowVals[0,5] = 1.054380329205709
$ws.Range("I3:N3").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.048810849987736
$rowVals[0,2] = 1.055963728759096
$rowVals[0,3] = 1.055957023272398
$rowVals[0,4] = 1.066409090969609
$ws.Range("B4:F4").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041934314458932
$rowVals[0,1] = 1.053281184469461
$rowVals[0,2] = 1.058398233924547
$rowVals[0,3] = 1.058391544712878
$rowVals[0,4] = 1.068818510833508
$rowVals[0,5] = 1.054776964564448
$ws.Range("I4:N4").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.049073337513307
$rowVals[0,2] = 1.056196462301687
$rowVals[0,3] = 1.056188478737683
$rowVals[0,4] = 1.06665838934288
$ws.Range("B5:F5").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.04197902752631
$rowVals[0,1] = 1.053447588087189
$rowVals[0,2] = 1.058579717085143
$rowVals[0,3] = 1.058571752467039
$rowVals[0,4] = 1.069017078167228
$rowVals[0,5] = 1.0549436044944
$ws.Range("I5:N5").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.049117412951672
$rowVals[0,2] = 1.056235545712303
$rowVals[0,3] = 1.056227348447263
$rowVals[0,4] = 1.066700254770582
$ws.Range("B6:F6").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041986508553277
$rowVals[0,1] = 1.053475521774921
$rowVals[0,2] = 1.058610186299473
$rowVals[0,3] = 1.058602008413155
$rowVals[0,4] = 1.06905041674371
$rowVals[0,5] = 1.054971577851173
$ws.Range("I6:N6").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.048814357182147
$rowVals[0,2] = 1.055966838124651
$rowVals[0,3] = 1.055960115500329
$rowVals[0,4] = 1.066412421627687
$ws.Range("B7:F7").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041934913691899
$rowVals[0,1] = 1.053283408378282
$rowVals[0,2] = 1.058400659085343
$rowVals[0,3] = 1.058393952772792
$rowVals[0,4] = 1.06882116421654
$rowVals[0,5] = 1.054779191631474
$ws.Range("I7:N7").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.047548020568561
$rowVals[0,2] = 1.054844623427797
$rowVals[0,3] = 1.054844193961761
$rowVals[0,4] = 1.065210368260329
$ws.Range("B8:F8").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041715426216682
$rowVals[0,1] = 1.052479513787404
$rowVals[0,2] = 1.057524499876415
$rowVals[0,3] = 1.057524071564523
$rowVals[0,4] = 1.067862680823009
$rowVals[0,5] = 1.053974155418069
$ws.Range("I8:N8").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.045318095953073
$rowVals[0,2] = 1.052870766366006
$rowVals[0,3] = 1.052881914531584
$rowVals[0,4] = 1.063096217405814
$ws.Range("B9:F9").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041314273059709
$rowVals[0,1] = 1.051059653394548
$rowVals[0,2] = 1.055979274175746
$rowVals[0,3] = 1.055990387238473
$rowVals[0,4] = 1.066172865579914
$rowVals[0,5] = 1.052552278660595
$ws.Range("I9:N9").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.043832558219635
$rowVals[0,2] = 1.05155737975778
$rowVals[0,3] = 1.051576575148722
$rowVals[0,4] = 1.061689568558229
$ws.Range("B10:F10").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.04103720512478
$rowVals[0,1] = 1.050110923134937
$rowVals[0,2] = 1.054948299240238
$rowVals[0,3] = 1.054967428265823
$rowVals[0,4] = 1.065045824232682
$rowVals[0,5] = 1.051602201095195
$ws.Range("I10:N10").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.043189570837474
$rowVals[0,2] = 1.050989282256328
$rowVals[0,3] = 1.051012040298942
$rowVals[0,4] = 1.061081151458018
$ws.Range("B11:F11").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.040914955379068
$rowVals[0,1] = 1.049699612968577
$rowVals[0,2] = 1.054501693238026
$rowVals[0,3] = 1.054524369551005
$rowVals[0,4] = 1.064557699221326
$rowVals[0,5] = 1.051190306821221
$ws.Range("I11:N11").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.042950776609083
$rowVals[0,2] = 1.050778358089694
$rowVals[0,3] = 1.050802451213908
$rowVals[0,4] = 1.060855260397048
$ws.Range("B12:F12").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.040869204840862
$rowVals[0,1] = 1.049546759304024
$rowVals[0,2] = 1.054335776789357
$rowVals[0,3] = 1.054359782190234
$rowVals[0,4] = 1.064376372645946
$rowVals[0,5] = 1.051037236086931
$ws.Range("I12:N12").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.043001996962423
$rowVals[0,2] = 1.050823597831629
$rowVals[0,3] = 1.050847404044198
$rowVals[0,4] = 1.060903710153411
$ws.Range("B13:F13").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.040879033931932
$rowVals[0,1] = 1.049579550302035
$rowVals[0,2] = 1.054371367640744
$rowVals[0,3] = 1.054395087428817
$rowVals[0,4] = 1.064415268482443
$rowVals[0,5] = 1.05107007365192
$ws.Range("I13:N13").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.043169831213122
$rowVals[0,2] = 1.050971845303861
$rowVals[0,3] = 1.05099471346488
$rowVals[0,4] = 1.061062477141274
$ws.Range("B14:F14").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.040911180588985
$rowVals[0,1] = 1.049686979556657
$rowVals[0,2] = 1.054487979079523
$rowVals[0,3] = 1.054510765021385
$rowVals[0,4] = 1.064542711013942
$rowVals[0,5] = 1.051177655468407
$ws.Range("I14:N14").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.043273244747562
$rowVals[0,2] = 1.051063197764335
$rowVals[0,3] = 1.05108548951876
$rowVals[0,4] = 1.06116031229608
$ws.Range("B15:F15").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.040930941966827
$rowVals[0,1] = 1.049753160387724
$rowVals[0,2] = 1.054559823675443
$rowVals[0,3] = 1.054582035763403
$rowVals[0,4] = 1.064621230581192
$rowVals[0,5] = 1.051243930283845
$ws.Range("I15:N15").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.043875236727373
$rowVals[0,2] = 1.051595095426352
$rowVals[0,3] = 1.051614056023474
$rowVals[0,4] = 1.061729961449298
$ws.Range("B16:F16").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041045270507265
$rowVals[0,1] = 1.050138209902277
$rowVals[0,2] = 1.05497793515492
$rowVals[0,3] = 1.054996830366821
$rowVals[0,4] = 1.065078217287722
$rowVals[0,5] = 1.051629526612875
$ws.Range("I16:N16").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.044252920645334
$rowVals[0,2] = 1.051928904331045
$rowVals[0,3] = 1.051945796313327
$rowVals[0,4] = 1.062087467920679
$ws.Range("B17:F17").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041116376474058
$rowVals[0,1] = 1.050379607085373
$rowVals[0,2] = 1.055240156073583
$rowVals[0,3] = 1.055256991131031
$rowVals[0,4] = 1.065364844669087
$rowVals[0,5] = 1.051871266607659
$ws.Range("I17:N17").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.044473242341579
$rowVals[0,2] = 1.052123667954738
$rowVals[0,3] = 1.052139360888448
$rowVals[0,4] = 1.062296060063983
$ws.Range("B18:F18").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041157631593645
$rowVals[0,1] = 1.05052036131659
$rowVals[0,2] = 1.055393086983865
$rowVals[0,3] = 1.05540872775204
$rowVals[0,4] = 1.065532019011689
$rowVals[0,5] = 1.05201222072603
$ws.Range("I18:N18").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.044548370571681
$rowVals[0,2] = 1.052190087207793
$rowVals[0,3] = 1.052205372575356
$rowVals[0,4] = 1.062367195523086
$ws.Range("B19:F19").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041171661224253
$rowVals[0,1] = 1.050568346590109
$rowVals[0,2] = 1.055445229368808
$rowVals[0,3] = 1.055460464164235
$rowVals[0,4] = 1.065589019317288
$rowVals[0,5] = 1.05206027414414
$ws.Range("I19:N19").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.044212396137608
$rowVals[0,2] = 1.051893083730781
$rowVals[0,3] = 1.051910196875924
$rowVals[0,4] = 1.062049104161374
$ws.Range("B20:F20").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.041108770205865
$rowVals[0,1] = 1.05035371247878
$rowVals[0,2] = 1.055212024117047
$rowVals[0,3] = 1.055229079447878
$rowVals[0,4] = 1.065334093347809
$rowVals[0,5] = 1.051845335227755
$ws.Range("I20:N20").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.043120407101091
$rowVals[0,2] = 1.050928187553608
$rowVals[0,3] = 1.050951331622852
$rowVals[0,4] = 1.061015721390456
$ws.Range("B21:F21").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.040901723622846
$rowVals[0,1] = 1.049655346375009
$rowVals[0,2] = 1.054453640657126
$rowVals[0,3] = 1.054476701276857
$rowVals[0,4] = 1.064505182775847
$rowVals[0,5] = 1.051145977364013
$ws.Range("I21:N21").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.042434061587149
$rowVals[0,2] = 1.05032205505842
$rowVals[0,3] = 1.050349059021293
$rowVals[0,4] = 1.060366584222104
$ws.Range("B22:F22").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.040769569484857
$rowVals[0,1] = 1.049215823830036
$rowVals[0,2] = 1.053976658612569
$rowVals[0,3] = 1.054003561343175
$rowVals[0,4] = 1.063983926334055
$rowVals[0,5] = 1.050705830646612
$ws.Range("I22:N22").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.042797884131033
$rowVals[0,2] = 1.050643326133775
$rowVals[0,3] = 1.050668277442882
$rowVals[0,4] = 1.060710647604699
$ws.Range("B23:F23").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.040839813961531
$rowVals[0,1] = 1.04944886369486
$rowVals[0,2] = 1.054229530295567
$rowVals[0,3] = 1.054254389961246
$rowVals[0,4] = 1.064260262156235
$rowVals[0,5] = 1.050939201454774
$ws.Range("I23:N23").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.044230707340882
$rowVals[0,2] = 1.051909269337393
$rowVals[0,3] = 1.051926282525267
$rowVals[0,4] = 1.062066438893148
$ws.Range("B24:F24").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,6
$rowVals[0,0] = 1.0411122078302
$rowVals[0,1] = 1.050365413287292
$rowVals[0,2] = 1.055224735793322
$rowVals[0,3] = 1.055241691569342
$rowVals[0,4] = 1.065347988578927
$rowVals[0,5] = 1.051857052652757
$ws.Range("I24:N24").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.02
$rowVals[0,1] = 1.045894398825198
$rowVals[0,2] = 1.053380617619916
$rowVals[0,3] = 1.053388714773849
$rowVals[0,4] = 1.063642291156528
$ws.Range("B25:F25").Value = $rowVals

$rowVals = New-Object 'object[,]' 1,5
$rowVals[0,0] = 1.041419681565519
$rowVals[0,1] = 1.051427105768394
$rowVals[0,2] = 1.056378902034078
$rowVals[0,3] = 1.056386974659262
$rowVals[0,4] = 1.05292025285896
$ws.Range("I25:M25").Value = $rowVals

